# Weekly update: insert a new price-record row for "Poroto verde" (Vega Modelo
# de Temuco) ahead of the existing history, pushing the older rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 68; all rows from 68 downward shift to 69+.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the latest week's record.
$ws.Cells.Item(68, 1).Value = 10
$ws.Cells.Item(68, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(68, 3).Value = "La Araucanía"
$ws.Cells.Item(68, 4).Value = 44524
$ws.Cells.Item(68, 5).Value = 9
$ws.Cells.Item(68, 6).Value = 100112031
$ws.Cells.Item(68, 7).Value = "Poroto verde"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 50
$ws.Cells.Item(68, 11).Value = 2000
$ws.Cells.Item(68, 12).Value = 2000
$ws.Cells.Item(68, 13).Value = 2000
$ws.Cells.Item(68, 14).Value = "$/kilo"
$ws.Cells.Item(68, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(68, 16).Value = 2000
$ws.Cells.Item(68, 17).Value = 1
$ws.Cells.Item(68, 18).Value = "Hortaliza"
